# A new daily price record (Cebollín, Terminal La Palmera de La Serena) was
# added as the new first entry of the data block (previously starting at
# row 37). Inserting a whole row at row 37 pushes every existing record
# down by one (old row 37 -> new row 38, ... old row 135 -> new row 136),
# which matches the observed diff exactly. We then populate the newly
# inserted row 37 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 37:135 down to 38:136, keeping their existing values/formats.
$ws.Rows.Item(37).Insert()

# Fill in the new record in row 37.
$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44526
$ws.Cells.Item(37, 5).Value = 4
$ws.Cells.Item(37, 6).Value = 100112037
$ws.Cells.Item(37, 7).Value = "Cebollín"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 3000
$ws.Cells.Item(37, 11).Value = 900
$ws.Cells.Item(37, 12).Value = 1000
$ws.Cells.Item(37, 13).Value = 950
$ws.Cells.Item(37, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(37, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 16).Value = 158
$ws.Cells.Item(37, 17).Value = 6
$ws.Cells.Item(37, 18).Value = "Hortaliza"
